$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value of 45182 (2023-09-13)
# for every data row (rows 2 through 265). Update it to 45184 (2023-09-15).
$ws.Range("C2:C265").Value = 45184
